$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.994.08"
$ws.Range("E2").Value = "  -0.48%  "

# Row 3
$ws.Range("D3").Value = "1.828.20"
$ws.Range("E3").Value = "  +0.08%  "

# Row 4
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").Value = "'311.57"
$ws.Range("E5").Value = "  -0.44%  "

# Row 6
$ws.Range("E6").Value = "  -0.16%  "

# Row 7
$ws.Range("D7").Value = "'0.4650"
$ws.Range("E7").Value = "  -1.22%  "

# Row 8
$ws.Range("D8").Value = "'0.3701"
$ws.Range("E8").Value = "  +1.53%  "

# Row 9
$ws.Range("D9").Value = "'0.07360"
$ws.Range("E9").Value = "  -0.46%  "

# Row 10
$ws.Range("D10").Value = "'0.8738"
$ws.Range("E10").Value = "  -0.68%  "

# Row 11
$ws.Range("D11").Value = "'0.07886"
$ws.Range("E11").Value = "  +7.49%  "

# Row 12
$ws.Range("D12").Value = "'19.93"
$ws.Range("E12").Value = "  -1.95%  "

# Row 13
$ws.Range("D13").Value = "1.862.38"
$ws.Range("E13").Value = "  -3.52%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.357"
$ws.Range("E14").Value = "  -0.41%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'6.579"
$ws.Range("E15").Value = "  +0.94%  "

# Row 16
$ws.Range("D16").Value = "'91.92"
$ws.Range("E16").Value = "  -1.46%  "

# Row 17
$ws.Range("D17").Value = "'1.009"
$ws.Range("E17").Value = "  +0.14%  "

# Row 18
$ws.Range("D18").Value = "'0.000008873"
$ws.Range("E18").Value = "  +1.97%  "

# Row 19
$ws.Range("D19").Value = "'1.007"
$ws.Range("E19").Value = "  -0.25%  "

# Row 20
$ws.Range("D20").Value = "'14.69"

# Row 21
$ws.Range("D21").Value = "26.917.93"
$ws.Range("E21").Value = "  -2.68%  "

# Row 22
$ws.Range("D22").Value = "'5.156"
$ws.Range("E22").Value = "  -1.64%  "

# Row 23
$ws.Range("D23").Value = "'10.57"
$ws.Range("E23").Value = "  -0.13%  "

# Row 24
$ws.Range("D24").Value = "2.087.05"
$ws.Range("E24").Value = "  -0.57%  "

# Row 25
$ws.Range("D25").Value = "'152.69"
$ws.Range("E25").Value = "  +0.64%  "

# Row 26
$ws.Range("D26").Value = "'1.830"
$ws.Range("E26").Value = "  -2.68%  "

# Row 27
$ws.Range("D27").Value = "'18.26"
$ws.Range("E27").Value = "  -1.44%  "

# Row 28
$ws.Range("D28").Value = "'2.102"
$ws.Range("E28").Value = "  -1.51%  "

# Row 29
$ws.Range("D29").Value = "'5.127"
$ws.Range("E29").Value = "  -0.98%  "

# Row 30
$ws.Range("D30").Value = "'115.44"

# Row 31
$ws.Range("D31").Value = "'0.08873"
$ws.Range("E31").Value = "  -0.72%  "

# Row 32
$ws.Range("D32").Value = "'2.982"
$ws.Range("E32").Value = "  +1.56%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7276"
$ws.Range("E33").Value = "  -1.77%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.443"
$ws.Range("E34").Value = "  -1.39%  "

# Row 35
$ws.Range("E35").Value = "  -2.73%  "

# Row 36
$ws.Range("D36").Value = "'2.506"
$ws.Range("E36").Value = "  +4.09%  "

# Row 37
$ws.Range("D37").Value = "'1.077"
$ws.Range("E37").Value = "  -1.00%  "

# Row 38
$ws.Range("E38").Value = "  +0.37%  "

# Row 39
$ws.Range("D39").Value = "'0.05235"
$ws.Range("E39").Value = "  -1.06%  "

# Row 40
$ws.Range("D40").Value = "'7.305"
$ws.Range("E40").Value = "  +1.71%  "

# Row 41
$ws.Range("D41").Value = "'2.928"
$ws.Range("E41").Value = "  -0.22%  "

# Row 42
$ws.Range("D42").Value = "'0.5185"
$ws.Range("E42").Value = "  -1.22%  "

# Row 43
$ws.Range("D43").Value = "'0.8604"
$ws.Range("E43").Value = "  -14.80%  "

# Row 44
$ws.Range("D44").Value = "'0.1625"
$ws.Range("E44").Value = "  -1.07%  "

# Row 45
$ws.Range("D45").Value = "'8.210"
$ws.Range("E45").Value = "  -1.89%  "

# Row 46
$ws.Range("D46").Value = "'0.4838"
$ws.Range("E46").Value = "  -0.46%  "

# Row 47
$ws.Range("D47").Value = "'1.008"
$ws.Range("E47").Value = "  -0.17%  "

# Row 48
$ws.Range("D48").Value = "'10.17"
$ws.Range("E48").Value = "  -1.63%  "

# Row 49
$ws.Range("D49").Value = "'102.75"
$ws.Range("E49").Value = "  -1.40%  "

# Row 50
$ws.Range("D50").Value = "'1.624"
$ws.Range("E50").Value = "  -1.54%  "

# Row 51
$ws.Range("D51").Value = "'0.06223"
$ws.Range("E51").Value = "  -1.17%  "
